$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C entirely - this shifts old D -> C, old E -> D
$ws.Range("C:C").Delete()

# Update header row
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"

# Update the B column (confidence-like values) for rows 2-14
$values = @(
    28724.4160595666,
    96498.87512547724,
    96484.65499409501,
    81949.7884529935,
    48389.93740341972,
    104549.8678707865,
    102228.8371356339,
    81280.10152616494,
    82870.82047733865,
    87521.31450254522,
    87480.36544053428,
    121875.6974668384,
    50698.79915833386
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
